$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.582.72'
$ws.Range("E2").Value = '  -1.02%  '

$ws.Range("D3").Value = '1.663.14'
$ws.Range("E3").Value = '  -3.68%  '

$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").Value = "'214.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.05%  '

$ws.Range("E6").Value = '  -1.97%  '

$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("D8").Value = "'23.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.44%  '

$ws.Range("E9").Value = '  -2.30%  '

$ws.Range("E10").Value = '  -1.59%  '

$ws.Range("E11").Value = '  -2.79%  '

$ws.Range("D13").Value = '1.661.85'
$ws.Range("E13").Value = '  -3.70%  '

$ws.Range("E14").Value = '  -2.86%  '

$ws.Range("E15").Value = '  -3.32%  '

$ws.Range("D16").Value = "'65.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.14%  '

$ws.Range("D17").Value = "'246.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.19%  '

$ws.Range("D18").Value = '27.566.13'
$ws.Range("E18").Value = '  -1.02%  '

$ws.Range("D19").Value = '0.0₃0733'
$ws.Range("E19").Value = '  -2.39%  '

$ws.Range("D20").Value = "'7.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.74%  '

$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("E22").Value = '  -3.80%  '

$ws.Range("E23").Value = '  -3.71%  '

$ws.Range("E24").Value = '  -4.93%  '

$ws.Range("D25").Value = "'146.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").Value = "'7.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.27%  '

$ws.Range("E27").Value = '  -2.71%  '

$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.20%  '

$ws.Range("D29").Value = "'0.112"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.16%  '

$ws.Range("E30").Value = '  +4.94%  '

$ws.Range("E31").Value = '  -1.34%  '

$ws.Range("E32").Value = '  -3.38%  '

$ws.Range("D33").Value = '1.449.80'
$ws.Range("E33").Value = '  -1.96%  '

$ws.Range("E34").Value = '  -5.23%  '

$ws.Range("E35").Value = '  -7.18%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = "'0.933"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.52%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = "'2.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.98%  '

$ws.Range("E38").Value = '  -5.82%  '

$ws.Range("E39").Value = '  -2.93%  '

$ws.Range("E40").Value = '  -2.96%  '

$ws.Range("D41").Value = "'69.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.96%  '

$ws.Range("E42").Value = '  +0.23%  '

$ws.Range("E43").Value = '  -8.18%  '

$ws.Range("D44").Value = "'0.791"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = "'2.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.32%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.805.20'
$ws.Range("E46").Value = '  -3.43%  '

$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").Value = "'88.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.03%  '

$ws.Range("D49").Value = '0.0₆0109'
$ws.Range("E49").Value = '  -1.11%  '

$ws.Range("E50").Value = '  -4.37%  '

$ws.Range("D51").Value = "'7.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.15%  '
